$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4287.5
$ws.Range("I69").Value = 4660
$ws.Range("J69").Value = 3666.6667
$ws.Range("K69").Value = 13980
$ws.Range("L69").Value = 11000.0001
$ws.Range("M69").Value = -13106
$ws.Range("N69").Value = -12748.0001

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4287.5
$ws.Range("I72").Value = 4660
$ws.Range("J72").Value = 3666.6667
$ws.Range("K72").Value = 41940
$ws.Range("L72").Value = 33000.0003
$ws.Range("M72").Value = -37572
$ws.Range("N72").Value = -41736.0003

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5949.9644
$ws.Range("I76").Value = 4605
$ws.Range("J76").Value = 9312.375
$ws.Range("K76").Value = 4605
$ws.Range("L76").Value = 9312.375
$ws.Range("M76").Value = -4290
$ws.Range("N76").Value = -9942.375

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5949.9644
$ws.Range("I79").Value = 4605
$ws.Range("J79").Value = 9312.375
$ws.Range("K79").Value = 4605
$ws.Range("L79").Value = 9312.375
$ws.Range("M79").Value = -3513
$ws.Range("N79").Value = -11496.375

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1481.4667
$ws.Range("I100").Value = 687.1429000000001
$ws.Range("J100").Value = 2176.5
$ws.Range("K100").Value = 687.1429000000001
$ws.Range("L100").Value = 2176.5
$ws.Range("M100").Value = -146.1429000000001
$ws.Range("N100").Value = -3258.5

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2783.8462
$ws.Range("I88").Value = 2723.75
$ws.Range("J88").Value = 2880
$ws.Range("K88").Value = 2723.75
$ws.Range("L88").Value = 2880
$ws.Range("M88").Value = -2317.75
$ws.Range("N88").Value = -3692

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2783.8462
$ws.Range("I91").Value = 2723.75
$ws.Range("J91").Value = 2880
$ws.Range("K91").Value = 2723.75
$ws.Range("L91").Value = 2880
$ws.Range("M91").Value = -1319.75
$ws.Range("N91").Value = -5688

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 5161.385
$ws.Range("I97").Value = 6256.9414
$ws.Range("J97").Value = 3092
$ws.Range("K97").Value = 6256.9414
$ws.Range("L97").Value = 3092
$ws.Range("M97").Value = -5760.9414
$ws.Range("N97").Value = -4084

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1815.8334
$ws.Range("I102").Value = 1498.5714
$ws.Range("J102").Value = 2260
$ws.Range("K102").Value = 1498.5714
$ws.Range("L102").Value = 2260
$ws.Range("M102").Value = 123.4286
$ws.Range("N102").Value = -5504

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8931378
$ws.Range("I132").Value = 17859758
$ws.Range("J132").Value = 2996.7144
$ws.Range("K132").Value = 53579274
$ws.Range("L132").Value = 8990.143199999999
$ws.Range("M132").Value = -53576744
$ws.Range("N132").Value = -14050.1432

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22729434
$ws.Range("I86").Value = 1809.091
$ws.Range("J86").Value = 45457060
$ws.Range("K86").Value = 1809.091
$ws.Range("L86").Value = 45457060
$ws.Range("M86").Value = -686.0909999999999
$ws.Range("N86").Value = -45459306

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 22729434
$ws.Range("I89").Value = 1809.091
$ws.Range("J89").Value = 45457060
$ws.Range("K89").Value = 9045.455
$ws.Range("L89").Value = 227285300
$ws.Range("M89").Value = -3429.455
$ws.Range("N89").Value = -227296532

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1787.4286
$ws.Range("I94").Value = 1622.75
$ws.Range("J94").Value = 2007
$ws.Range("K94").Value = 1622.75
$ws.Range("L94").Value = 2007
$ws.Range("M94").Value = -1171.75
$ws.Range("N94").Value = -2909

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4404.1113
$ws.Range("I105").Value = 2516.6667
$ws.Range("J105").Value = 4943.381
$ws.Range("K105").Value = 2516.6667
$ws.Range("L105").Value = 4943.381
$ws.Range("M105").Value = -769.6667000000002
$ws.Range("N105").Value = -8437.381000000001

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3546.3547
$ws.Range("I134").Value = 2326.5417
$ws.Range("K134").Value = 6979.625100000001
$ws.Range("M134").Value = -4444.625100000001

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2020
$ws.Range("I105").Value = 2275
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 2275
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -528
$ws.Range("N105").Value = -4494

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1823.55
$ws.Range("I134").Value = 2125.6428
$ws.Range("J134").Value = 1118.6666
$ws.Range("K134").Value = 6376.928400000001
$ws.Range("L134").Value = 3355.9998
$ws.Range("M134").Value = -3841.928400000001
$ws.Range("N134").Value = -8425.9998

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 37510.57
$ws.Range("J140").Value = 37510.57
$ws.Range("L140").Value = 37510.57
$ws.Range("N140").Value = -47870.57

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2996.1538
$ws.Range("I20").Value = 2900
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 8700
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -8473
$ws.Range("N20").Value = -9454

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 9944.444
$ws.Range("J106").Value = 9944.444
$ws.Range("L106").Value = 29833.332
$ws.Range("N106").Value = -31725.332

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1427.8182
$ws.Range("I113").Value = 575.8
$ws.Range("J113").Value = 2137.8333
$ws.Range("K113").Value = 1727.4
$ws.Range("L113").Value = 6413.499899999999
$ws.Range("M113").Value = 442.6000000000001
$ws.Range("N113").Value = -10753.4999

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 4019.9
$ws.Range("J125").Value = 4022.111
$ws.Range("L125").Value = 12066.333
$ws.Range("N125").Value = -21906.333

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2243.2285
$ws.Range("I134").Value = 1071.8928
$ws.Range("J134").Value = 6928.5713
$ws.Range("K134").Value = 3215.6784
$ws.Range("L134").Value = 20785.7139
$ws.Range("M134").Value = 1854.3216
$ws.Range("N134").Value = -30925.7139

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4138.52
$ws.Range("J137").Value = 6919.4165
$ws.Range("L137").Value = 20758.2495
$ws.Range("N137").Value = -30958.2495

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34050
$ws.Range("I70").Value = 201500
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 201500
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -201230
$ws.Range("N70").Value = -5040

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 34050
$ws.Range("I73").Value = 201500
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 201500
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -200564
$ws.Range("N73").Value = -6372

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 38195
$ws.Range("J88").Value = 38195
$ws.Range("L88").Value = 38195
$ws.Range("N88").Value = -39097

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 38195
$ws.Range("J91").Value = 38195
$ws.Range("L91").Value = 38195
$ws.Range("N91").Value = -41315

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1004.8461
$ws.Range("I16").Value = 1092.091
$ws.Range("K16").Value = 1092.091
$ws.Range("M16").Value = -922.0909999999999

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1685.7142
$ws.Range("I93").Value = 1475
$ws.Range("J93").Value = 1966.6666
$ws.Range("K93").Value = 1475
$ws.Range("L93").Value = 1966.6666
$ws.Range("M93").Value = -227
$ws.Range("N93").Value = -4462.6666

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6767.222
$ws.Range("I122").Value = 7901.2
$ws.Range("J122").Value = 5349.75
$ws.Range("K122").Value = 23703.6
$ws.Range("L122").Value = 16049.25
$ws.Range("M122").Value = -21253.6
$ws.Range("N122").Value = -20949.25

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10007020
$ws.Range("I132").Value = 4179.207
$ws.Range("J132").Value = 23820468
$ws.Range("K132").Value = 12537.621
$ws.Range("L132").Value = 71461404
$ws.Range("M132").Value = -10007.621
$ws.Range("N132").Value = -71466464

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11632508
$ws.Range("I136").Value = 17859002
$ws.Range("J136").Value = 9720.666999999999
$ws.Range("K136").Value = 53577006
$ws.Range("L136").Value = 29162.001
$ws.Range("M136").Value = -53574456
$ws.Range("N136").Value = -34262.001

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1124.9744
$ws.Range("I136").Value = 956.7879
$ws.Range("K136").Value = 2870.3637
$ws.Range("M136").Value = -320.3636999999999
